$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the A56 timestamp's tiny floating-point drift.
$ws.Cells.Item(56, 1).Value = 44369.76721021759

# Append the new day's row of job numbers.
$row = 57
$ws.Cells.Item($row, 1).Value = 44370.7665353545
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
$ws.Cells.Item($row, 2).Value = 78885
$ws.Cells.Item($row, 3).Value = 66307
$ws.Cells.Item($row, 4).Value = 3661
$ws.Cells.Item($row, 5).Value = 2127
$ws.Cells.Item($row, 6).Value = 1513
$ws.Cells.Item($row, 7).Value = 20800
$ws.Cells.Item($row, 8).Value = 1556
$ws.Cells.Item($row, 9).Value = 896
$ws.Cells.Item($row, 10).Value = 197
